$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "data2" worksheet right after the existing "data1" sheet.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "data2"

# Header row
$ws2.Range("A1").Value = "Địa điểm "
$ws2.Range("B1").Value = "Phòng"
$ws2.Range("C1").Value = "Nhận phòng"
$ws2.Range("D1").Value = "Trả phòng"
$ws2.Range("E1").Value = "Số khách"

# Data row
$ws2.Range("A2").Value = "Đà Lạt"
$ws2.Range("B2").Value = "Phòng mùa hè"
$ws2.Range("C2").Value = (Get-Date -Year 2025 -Month 10 -Day 1).Date
$ws2.Range("D2").Value = (Get-Date -Year 2025 -Month 10 -Day 15).Date
$ws2.Range("E2").Value = 2

# Column widths to roughly match the source layout
# (ColumnWidth values are pre-compensated for the host's implicit
#  "+5/7 char" padding that is added when the width is serialized.)
$ws2.Columns.Item(2).ColumnWidth = 13.082589285714286
$ws2.Columns.Item(3).ColumnWidth = 14.285714285714286
$ws2.Columns.Item(4).ColumnWidth = 15.582589285714286

# Select E3 on the new sheet (becomes the active/visible tab)
$ws2.Range("E3").Select()
